# Adds a new data column for the 12. 10. 2021 survey wave to both sheets
# (column AJ on "data", column AI on "pocetR"), and bumps the "aktualizace"
# (last-updated) date in each sheet's trailing footnote row from
# 6. 10. 2021 to 20. 10. 2021.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": percentages impacted, new column AJ (col 36), rows 1-67
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# New header cell AJ1 inherits its formatting (bold, centered, bordered)
# from the previous wave's header cell AI1.
$wsData.Range("AI1").Copy()
$wsData.Range("AJ1").PasteSpecial(-4122)
$wsData.Range("AJ1").Value() = "12. 10. 2021"

$dataNewCol = @{
    2 = 0.09
    3 = 0.17
    4 = 0.74
    5 = 0.07000000000000001
    6 = 0.15
    7 = 0.78
    8 = 0.05
    9 = 0.09
    10 = 0.86
    11 = 0.09
    12 = 0.22
    13 = 0.6899999999999999
    14 = 0.23
    15 = 0.15
    16 = 0.62
    17 = 0.06
    18 = 0.15
    19 = 0.79
    20 = 0.06
    21 = 0.08
    22 = 0.86
    23 = 0.19
    24 = 0.18
    25 = 0.63
    26 = 0.11
    27 = 0.22
    28 = 0.67
    29 = 0.15
    30 = 0.27
    31 = 0.58
    32 = 0.06
    33 = 0.15
    34 = 0.79
    35 = 0.03
    36 = 0.05
    37 = 0.92
    38 = 0.16
    39 = 0.2
    40 = 0.64
    41 = 0.07000000000000001
    42 = 0.17
    43 = 0.76
    44 = 0.65
    45 = 0.16
    46 = 0.19
    47 = 0.1
    48 = 0.5600000000000001
    49 = 0.34
    50 = 0.02
    51 = 0.08
    52 = 0.9
    53 = 0.06
    54 = 0.19
    55 = 0.75
    56 = 0.03
    57 = 0.07000000000000001
    58 = 0.9
    59 = 0.07000000000000001
    60 = 0.19
    61 = 0.74
    62 = 0.05
    63 = 0.07000000000000001
    64 = 0.88
    65 = 0.07000000000000001
    66 = 0.11
    67 = 0.82
}

foreach ($row in $dataNewCol.Keys) {
    $wsData.Cells.Item($row, 36).Value() = $dataNewCol[$row]
}

$wsData.Range("A68").Value() = "Život během pandemie, Zasažení domácností, % respondentů celkově a ve skupinách, aktualizace 20. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": sample sizes, new column AI (col 35), rows 1-23
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

# New header cell AI1 inherits its formatting from the previous wave's
# header cell AH1.
$wsPocet.Range("AH1").Copy()
$wsPocet.Range("AI1").PasteSpecial(-4122)
$wsPocet.Range("AI1").Value() = "12. 10. 2021"

$pocetNewCol = @{
    2 = 1562
    3 = 741
    4 = 123
    5 = 491
    6 = 207
    7 = 706
    8 = 114
    9 = 103
    10 = 639
    11 = 745
    12 = 514
    13 = 303
    14 = 426
    15 = 1136
    16 = 144
    17 = 280
    18 = 1138
    19 = 279
    20 = 93
    21 = 239
    22 = 143
    23 = 76
}

foreach ($row in $pocetNewCol.Keys) {
    $wsPocet.Cells.Item($row, 35).Value() = $pocetNewCol[$row]
}

$wsPocet.Range("A24").Value() = "Život během pandemie, Zasažení domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 20. 10. 2021"

# Row 24 is a footer row whose other trailing cells (B24:AH24) are present
# but empty, just to pad out the row. Mirror that for the new AI24 cell by
# copying the (empty) formatting of its left neighbour, AH24.
$wsPocet.Range("AH24").Copy()
$wsPocet.Range("AI24").PasteSpecial(-4122)
